# Regenerate penyata (statement) content to follow the new wording/data.
#
# Summary of the change (per the authoritative diff):
#   1. Label wording was made more specific: "Kali Pertama/Kedua/Ketiga/Keempat"
#      -> "Semakan Kali Pertama/Kedua/Ketiga/Keempat" everywhere those labels
#      are used (Merit Pendahuluan, JPPM/JDM/JDRM and Penyertaan Pertandingan
#      sections all reuse the same four labels).
#   2. Competition names were re-cased from ALL CAPS to Title Case.
#   3. A handful of transaction figures were updated with the latest data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Kali ..." -> "Semakan Kali ..." (every cell that carries the label)
# ---------------------------------------------------------------------
$kaliMap = @{
    "C16" = "Semakan Kali Pertama"
    "C22" = "Semakan Kali Pertama"
    "C28" = "Semakan Kali Pertama"

    "C17" = "Semakan Kali Kedua"
    "C23" = "Semakan Kali Kedua"
    "C29" = "Semakan Kali Kedua"

    "C18" = "Semakan Kali Ketiga"
    "C24" = "Semakan Kali Ketiga"
    "C30" = "Semakan Kali Ketiga"

    "C19" = "Semakan Kali Keempat"
    "C25" = "Semakan Kali Keempat"
    "C31" = "Semakan Kali Keempat"
}
foreach ($addr in $kaliMap.Keys) {
    $ws.Range($addr).Value = $kaliMap[$addr]
}

# ---------------------------------------------------------------------
# 2) Competition names: ALL CAPS -> Title Case
# ---------------------------------------------------------------------
$ws.Range("C34").Value = "Bouquet Kreatif"
$ws.Range("C35").Value = "Tik Tok Raya"
$ws.Range("C36").Value = "Riang Ria Kuih Raya"
$ws.Range("C37").Value = "Creative Collage"

# ---------------------------------------------------------------------
# 3) Updated transaction figures
# ---------------------------------------------------------------------
$ws.Range("D18").Value = 10235
$ws.Range("E18").Value = 950
$ws.Range("E23").Value = 200
$ws.Range("E29").Value = 1750
